$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 8.841467
$ws.Range("H2").Value = 26.524401
$ws.Range("I2").Value = 0.5917001192060068
$ws.Range("J2").Value = 0.5917001192060067
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 19.11595033333333
$ws.Range("N2").Value = 57.347851
$ws.Range("O2").Value = 0.6851940154453416
$ws.Range("P2").Value = 0.6851940154453418
$ws.Range("Q2").Value = 169.0130440458057
$ws.Range("R2").Value = 1521.117396412251
$ws.Range("S2").Value = 0.4054293806182511
$ws.Range("T2").Value = 0.4054293806182511

# Row 3
$ws.Range("G3").Value = 8.841467
$ws.Range("H3").Value = 26.524401
$ws.Range("I3").Value = 0.5917001192060068
$ws.Range("J3").Value = 0.5917001192060067
$ws.Range("M3").Value = 4.865208333333334
$ws.Range("O3").Value = 0.1743890089566637
$ws.Range("P3").Value = 0.1743890089566637
$ws.Range("Q3").Value = 43.01557892729167
$ws.Range("R3").Value = 387.1402103456251
$ws.Range("S3").Value = 0.1031859973878753
$ws.Range("T3").Value = 0.1031859973878753

# Row 4
$ws.Range("G4").Value = 8.841467
$ws.Range("H4").Value = 26.524401
$ws.Range("I4").Value = 0.5917001192060068
$ws.Range("J4").Value = 0.5917001192060067
$ws.Range("M4").Value = 3.917436333333333
$ws.Range("N4").Value = 11.752309
$ws.Range("O4").Value = 0.1404169755979945
$ws.Range("P4").Value = 0.1404169755979946
$ws.Range("Q4").Value = 34.63588406576766
$ws.Range("R4").Value = 311.722956591909
$ws.Range("S4").Value = 0.08308474119988032
$ws.Range("T4").Value = 0.08308474119988032

# Row 5
$ws.Range("I5").Value = 0.2746155987184545
$ws.Range("J5").Value = 0.2746155987184545
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 19.11595033333333
$ws.Range("N5").Value = 57.347851
$ws.Range("O5").Value = 0.6851940154453416
$ws.Range("P5").Value = 0.6851940154453418
$ws.Range("Q5").Value = 78.44111700391267
$ws.Range("R5").Value = 705.9700530352138
$ws.Range("S5").Value = 0.1881649647898245
$ws.Range("T5").Value = 0.1881649647898245

# Row 6
$ws.Range("I6").Value = 0.2746155987184545
$ws.Range("J6").Value = 0.2746155987184545
$ws.Range("M6").Value = 4.865208333333334
$ws.Range("O6").Value = 0.1743890089566637
$ws.Range("P6").Value = 0.1743890089566637
$ws.Range("S6").Value = 0.04788994210455214
$ws.Range("T6").Value = 0.04788994210455213

# Row 7
$ws.Range("I7").Value = 0.2746155987184545
$ws.Range("J7").Value = 0.2746155987184545
$ws.Range("M7").Value = 3.917436333333333
$ws.Range("N7").Value = 11.752309
$ws.Range("O7").Value = 0.1404169755979945
$ws.Range("P7").Value = 0.1404169755979946
$ws.Range("Q7").Value = 16.07495711278067
$ws.Range("R7").Value = 144.674614015026
$ws.Range("S7").Value = 0.03856069182407789
$ws.Range("T7").Value = 0.03856069182407789

# Row 8
$ws.Range("G8").Value = 1.997574666666667
$ws.Range("H8").Value = 5.992724
$ws.Range("I8").Value = 0.1336842820755386
$ws.Range("J8").Value = 0.1336842820755386
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 19.11595033333333
$ws.Range("N8").Value = 57.347851
$ws.Range("O8").Value = 0.6851940154453416
$ws.Range("P8").Value = 0.6851940154453418
$ws.Range("Q8").Value = 38.18553811512489
$ws.Range("R8").Value = 343.669843036124
$ws.Range("S8").Value = 0.09159967003726602
$ws.Range("T8").Value = 0.09159967003726602

# Row 9
$ws.Range("G9").Value = 1.997574666666667
$ws.Range("H9").Value = 5.992724
$ws.Range("I9").Value = 0.1336842820755386
$ws.Range("J9").Value = 0.1336842820755386
$ws.Range("M9").Value = 4.865208333333334
$ws.Range("O9").Value = 0.1743890089566637
$ws.Range("P9").Value = 0.1743890089566637
$ws.Range("Q9").Value = 9.718616914722224
$ws.Range("R9").Value = 87.46755223250001
$ws.Range("S9").Value = 0.02331306946423626
$ws.Range("T9").Value = 0.02331306946423626

# Row 10
$ws.Range("G10").Value = 1.997574666666667
$ws.Range("H10").Value = 5.992724
$ws.Range("I10").Value = 0.1336842820755386
$ws.Range("J10").Value = 0.1336842820755386
$ws.Range("M10").Value = 3.917436333333333
$ws.Range("N10").Value = 11.752309
$ws.Range("O10").Value = 0.1404169755979945
$ws.Range("P10").Value = 0.1404169755979946
$ws.Range("Q10").Value = 7.825371577746222
$ws.Range("R10").Value = 70.42834419971601
$ws.Range("S10").Value = 0.01877154257403632
$ws.Range("T10").Value = 0.01877154257403632
